$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '305.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.20%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '35.44'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-2.52%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.101'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.77%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08020'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '1.11%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.974'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-10.10%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.137'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-0.82%'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '7.845'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-2.23%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.896'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '10.01%'
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9264'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-0.43%'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1092'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '10.74%'
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1895'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '0.78%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09421'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '3.80%'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03672'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.67%'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09971'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.53%'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001410'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.67%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.005833'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '3.57%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.457'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.59%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3416'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.33%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1311'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.54%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.129'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2196'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '0.31%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04529'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.49%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001232'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.58%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004679'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.15%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001260'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-3.26%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0004480'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-5.45%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01876'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-4.97%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04755'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-3.55%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007585'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-2.33%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.009674'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '23.92%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1346'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-3.61%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002131'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '0.78%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01132'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.95%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006390'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '2.98%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000756'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.74%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '64.42'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '24.11%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.001309'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-27.25%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002116'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.74%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002015'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.74%'
